# Update scripts with new TPM values: refresh ligand/receptor expression
# columns (G,H,M,N) and the derived specificity/edge-weight columns
# (I,J,O,P,Q,R,S,T) that cascade from them for every sending/target
# cluster combination on the L1cam-Egfr sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 7.741029
$ws.Cells.Item(2, 8).Value = 23.223087
$ws.Cells.Item(2, 9).Value = 0.4930486933812723
$ws.Cells.Item(2, 10).Value = 0.4930486933812723
$ws.Cells.Item(2, 13).Value = 0.4102596666666667
$ws.Cells.Item(2, 14).Value = 1.230779
$ws.Cells.Item(2, 15).Value = 0.003499619873322347
$ws.Cells.Item(2, 16).Value = 0.003499619873322347
$ws.Cells.Item(2, 17).Value = 3.175831977197
$ws.Cells.Item(2, 18).Value = 28.582487794773
$ws.Cells.Item(2, 19).Value = 0.001725483005872717
$ws.Cells.Item(2, 20).Value = 0.001725483005872717
$ws.Cells.Item(3, 7).Value = 7.741029
$ws.Cells.Item(3, 8).Value = 23.223087
$ws.Cells.Item(3, 9).Value = 0.4930486933812723
$ws.Cells.Item(3, 10).Value = 0.4930486933812723
$ws.Cells.Item(3, 15).Value = 0.8692174743460166
$ws.Cells.Item(3, 16).Value = 0.8692174743460165
$ws.Cells.Item(3, 17).Value = 788.796712240017
$ws.Cells.Item(3, 18).Value = 7099.170410160154
$ws.Cells.Item(3, 19).Value = 0.4285665399904731
$ws.Cells.Item(3, 20).Value = 0.428566539990473
$ws.Cells.Item(4, 7).Value = 7.741029
$ws.Cells.Item(4, 8).Value = 23.223087
$ws.Cells.Item(4, 9).Value = 0.4930486933812723
$ws.Cells.Item(4, 10).Value = 0.4930486933812723
$ws.Cells.Item(4, 14).Value = 44.764041
$ws.Cells.Item(4, 15).Value = 0.1272829057806611
$ws.Cells.Item(4, 16).Value = 0.1272829057806611
$ws.Cells.Item(4, 17).Value = 115.506579846063
$ws.Cells.Item(4, 18).Value = 1039.559218614567
$ws.Cells.Item(4, 19).Value = 0.06275667038492656
$ws.Cells.Item(4, 20).Value = 0.06275667038492656
$ws.Cells.Item(5, 9).Value = 0.0194007766416684
$ws.Cells.Item(5, 10).Value = 0.0194007766416684
$ws.Cells.Item(5, 13).Value = 0.4102596666666667
$ws.Cells.Item(5, 14).Value = 1.230779
$ws.Cells.Item(5, 15).Value = 0.003499619873322347
$ws.Cells.Item(5, 16).Value = 0.003499619873322347
$ws.Cells.Item(5, 17).Value = 0.1249645474537778
$ws.Cells.Item(5, 18).Value = 1.124680927084
$ws.Cells.Item(5, 19).Value = 0.00006789534349307072
$ws.Cells.Item(5, 20).Value = 0.00006789534349307071
$ws.Cells.Item(6, 9).Value = 0.0194007766416684
$ws.Cells.Item(6, 10).Value = 0.0194007766416684
$ws.Cells.Item(6, 15).Value = 0.8692174743460166
$ws.Cells.Item(6, 16).Value = 0.8692174743460165
$ws.Cells.Item(6, 19).Value = 0.0168634940728222
$ws.Cells.Item(6, 20).Value = 0.0168634940728222
$ws.Cells.Item(7, 9).Value = 0.0194007766416684
$ws.Cells.Item(7, 10).Value = 0.0194007766416684
$ws.Cells.Item(7, 14).Value = 44.764041
$ws.Cells.Item(7, 15).Value = 0.1272829057806611
$ws.Cells.Item(7, 16).Value = 0.1272829057806611
$ws.Cells.Item(7, 17).Value = 4.545022401070667
$ws.Cells.Item(7, 18).Value = 40.905201609636
$ws.Cells.Item(7, 19).Value = 0.002469387225353131
$ws.Cells.Item(7, 20).Value = 0.002469387225353131
$ws.Cells.Item(8, 7).Value = 7.654706000000001
$ws.Cells.Item(8, 9).Value = 0.4875505299770593
$ws.Cells.Item(8, 10).Value = 0.4875505299770593
$ws.Cells.Item(8, 13).Value = 0.4102596666666667
$ws.Cells.Item(8, 14).Value = 1.230779
$ws.Cells.Item(8, 15).Value = 0.003499619873322347
$ws.Cells.Item(8, 16).Value = 0.003499619873322347
$ws.Cells.Item(8, 17).Value = 3.140417131991334
$ws.Cells.Item(8, 18).Value = 28.263754187922
$ws.Cells.Item(8, 19).Value = 0.00170624152395656
$ws.Cells.Item(8, 20).Value = 0.00170624152395656
$ws.Cells.Item(9, 7).Value = 7.654706000000001
$ws.Cells.Item(9, 9).Value = 0.4875505299770593
$ws.Cells.Item(9, 10).Value = 0.4875505299770593
$ws.Cells.Item(9, 15).Value = 0.8692174743460166
$ws.Cells.Item(9, 16).Value = 0.8692174743460165
$ws.Cells.Item(9, 18).Value = 7020.005006269243
$ws.Cells.Item(9, 19).Value = 0.4237874402827214
$ws.Cells.Item(9, 20).Value = 0.4237874402827213
$ws.Cells.Item(10, 7).Value = 7.654706000000001
$ws.Cells.Item(10, 9).Value = 0.4875505299770593
$ws.Cells.Item(10, 10).Value = 0.4875505299770593
$ws.Cells.Item(10, 14).Value = 44.764041
$ws.Cells.Item(10, 15).Value = 0.1272829057806611
$ws.Cells.Item(10, 16).Value = 0.1272829057806611
$ws.Cells.Item(10, 19).Value = 0.06205684817038145
$ws.Cells.Item(10, 20).Value = 0.06205684817038145
